$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E5").Value = "❌ EXPIRED 3754 days ago"
$ws.Range("E6").Value = "⚠️ Expires in 25 days"
$ws.Range("E7").Value = "⚠️ Expires in 25 days"
$ws.Range("E8").Value = "⚠️ Expires in 25 days"
